$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recalculated TPM-based ligand/receptor/edge expression values
# (columns G,H,I,J = ligand; M,N,O,P = receptor; Q,R,S,T = edge)
    $ws.Range("G2").Value = 27.67634766666667
    $ws.Range("H2").Value = 83.029043
    $ws.Range("I2").Value = 0.005965811625935536
    $ws.Range("J2").Value = 0.005965811625935536
    $ws.Range("M2").Value = 3.087329333333333
    $ws.Range("N2").Value = 9.261987999999999
    $ws.Range("O2").Value = 0.1539049749041678
    $ws.Range("P2").Value = 0.1539049749041678
    $ws.Range("Q2").Value = 85.44599999083155
    $ws.Range("R2").Value = 769.0139999174839
    $ws.Range("S2").Value = 0.0009181680885726009
    $ws.Range("T2").Value = 0.0009181680885726009
    $ws.Range("G3").Value = 27.67634766666667
    $ws.Range("H3").Value = 83.029043
    $ws.Range("I3").Value = 0.005965811625935536
    $ws.Range("J3").Value = 0.005965811625935536
    $ws.Range("O3").Value = 0.2832552948356705
    $ws.Range("P3").Value = 0.2832552948356705
    $ws.Range("Q3").Value = 157.259581342333
    $ws.Range("R3").Value = 1415.336232080997
    $ws.Range("S3").Value = 0.001689847731038441
    $ws.Range("T3").Value = 0.001689847731038441
    $ws.Range("G4").Value = 27.67634766666667
    $ws.Range("H4").Value = 83.029043
    $ws.Range("I4").Value = 0.005965811625935536
    $ws.Range("J4").Value = 0.005965811625935536
    $ws.Range("M4").Value = 4.823431
    $ws.Range("N4").Value = 14.470293
    $ws.Range("O4").Value = 0.2404505470122564
    $ws.Range("P4").Value = 0.2404505470122564
    $ws.Range("Q4").Value = 133.4949533021777
    $ws.Range("R4").Value = 1201.454579719599
    $ws.Range("S4").Value = 0.001434482668828279
    $ws.Range("T4").Value = 0.001434482668828278
    $ws.Range("G5").Value = 27.67634766666667
    $ws.Range("H5").Value = 83.029043
    $ws.Range("I5").Value = 0.005965811625935536
    $ws.Range("J5").Value = 0.005965811625935536
    $ws.Range("M5").Value = 6.467117666666667
    $ws.Range("N5").Value = 19.401353
    $ws.Range("O5").Value = 0.3223891832479054
    $ws.Range("P5").Value = 0.3223891832479053
    $ws.Range("Q5").Value = 178.9861969439088
    $ws.Range("R5").Value = 1610.875772495179
    $ws.Range("S5").Value = 0.001923313137496216
    $ws.Range("T5").Value = 0.001923313137496215
    $ws.Range("I6").Value = 0.009118181457976757
    $ws.Range("J6").Value = 0.009118181457976757
    $ws.Range("M6").Value = 3.087329333333333
    $ws.Range("N6").Value = 9.261987999999999
    $ws.Range("O6").Value = 0.1539049749041678
    $ws.Range("P6").Value = 0.1539049749041678
    $ws.Range("Q6").Value = 130.5961672318986
    $ws.Range("R6").Value = 1175.365505087088
    $ws.Range("S6").Value = 0.001403333488461561
    $ws.Range("T6").Value = 0.001403333488461561
    $ws.Range("I7").Value = 0.009118181457976757
    $ws.Range("J7").Value = 0.009118181457976757
    $ws.Range("O7").Value = 0.2832552948356705
    $ws.Range("P7").Value = 0.2832552948356705
    $ws.Range("S7").Value = 0.00258277317724435
    $ws.Range("T7").Value = 0.00258277317724435
    $ws.Range("I8").Value = 0.009118181457976757
    $ws.Range("J8").Value = 0.009118181457976757
    $ws.Range("M8").Value = 4.823431
    $ws.Range("N8").Value = 14.470293
    $ws.Range("O8").Value = 0.2404505470122564
    $ws.Range("P8").Value = 0.2404505470122564
    $ws.Range("Q8").Value = 204.034469114252
    $ws.Range("R8").Value = 1836.310222028268
    $ws.Range("S8").Value = 0.002192471719327525
    $ws.Range("T8").Value = 0.002192471719327525
    $ws.Range("I9").Value = 0.009118181457976757
    $ws.Range("J9").Value = 0.009118181457976757
    $ws.Range("M9").Value = 6.467117666666667
    $ws.Range("N9").Value = 19.401353
    $ws.Range("O9").Value = 0.3223891832479054
    $ws.Range("P9").Value = 0.3223891832479053
    $ws.Range("Q9").Value = 273.5635525454254
    $ws.Range("R9").Value = 2462.071972908828
    $ws.Range("S9").Value = 0.002939603072943322
    $ws.Range("T9").Value = 0.002939603072943321
    $ws.Range("G10").Value = 29.593002
    $ws.Range("H10").Value = 88.779006
    $ws.Range("I10").Value = 0.006378958578792732
    $ws.Range("J10").Value = 0.006378958578792732
    $ws.Range("M10").Value = 3.087329333333333
    $ws.Range("N10").Value = 9.261987999999999
    $ws.Range("O10").Value = 0.1539049749041678
    $ws.Range("P10").Value = 0.1539049749041678
    $ws.Range("Q10").Value = 91.36334313599198
    $ws.Range("R10").Value = 822.2700882239278
    $ws.Range("S10").Value = 0.0009817534599838211
    $ws.Range("T10").Value = 0.0009817534599838211
    $ws.Range("G11").Value = 29.593002
    $ws.Range("H11").Value = 88.779006
    $ws.Range("I11").Value = 0.006378958578792732
    $ws.Range("J11").Value = 0.006378958578792732
    $ws.Range("O11").Value = 0.2832552948356705
    $ws.Range("P11").Value = 0.2832552948356705
    $ws.Range("Q11").Value = 168.150189513186
    $ws.Range("R11").Value = 1513.351705618674
    $ws.Range("S11").Value = 0.001806873792980465
    $ws.Range("T11").Value = 0.001806873792980465
    $ws.Range("G12").Value = 29.593002
    $ws.Range("H12").Value = 88.779006
    $ws.Range("I12").Value = 0.006378958578792732
    $ws.Range("J12").Value = 0.006378958578792732
    $ws.Range("M12").Value = 4.823431
    $ws.Range("N12").Value = 14.470293
    $ws.Range("O12").Value = 0.2404505470122564
    $ws.Range("P12").Value = 0.2404505470122564
    $ws.Range("Q12").Value = 142.739803229862
    $ws.Range("R12").Value = 1284.658229068758
    $ws.Range("S12").Value = 0.001533824079639238
    $ws.Range("T12").Value = 0.001533824079639238
    $ws.Range("G13").Value = 29.593002
    $ws.Range("H13").Value = 88.779006
    $ws.Range("I13").Value = 0.006378958578792732
    $ws.Range("J13").Value = 0.006378958578792732
    $ws.Range("M13").Value = 6.467117666666667
    $ws.Range("N13").Value = 19.401353
    $ws.Range("O13").Value = 0.3223891832479054
    $ws.Range("P13").Value = 0.3223891832479053
    $ws.Range("Q13").Value = 191.381426043902
    $ws.Range("R13").Value = 1722.432834395118
    $ws.Range("S13").Value = 0.002056507246189208
    $ws.Range("T13").Value = 0.002056507246189208
    $ws.Range("G14").Value = 4539.588785666667
    $ws.Range("H14").Value = 13618.766357
    $ws.Range("I14").Value = 0.9785370483372949
    $ws.Range("J14").Value = 0.978537048337295
    $ws.Range("M14").Value = 3.087329333333333
    $ws.Range("N14").Value = 9.261987999999999
    $ws.Range("O14").Value = 0.1539049749041678
    $ws.Range("P14").Value = 0.1539049749041678
    $ws.Range("Q14").Value = 14015.20561925974
    $ws.Range("R14").Value = 126136.8505733377
    $ws.Range("S14").Value = 0.1506017198671498
    $ws.Range("T14").Value = 0.1506017198671498
    $ws.Range("G15").Value = 4539.588785666667
    $ws.Range("H15").Value = 13618.766357
    $ws.Range("I15").Value = 0.9785370483372949
    $ws.Range("J15").Value = 0.978537048337295
    $ws.Range("O15").Value = 0.2832552948356705
    $ws.Range("P15").Value = 0.2832552948356705
    $ws.Range("Q15").Value = 25794.36566191506
    $ws.Range("R15").Value = 232149.2909572356
    $ws.Range("S15").Value = 0.2771758001344072
    $ws.Range("T15").Value = 0.2771758001344072
    $ws.Range("G16").Value = 4539.588785666667
    $ws.Range("H16").Value = 13618.766357
    $ws.Range("I16").Value = 0.9785370483372949
    $ws.Range("J16").Value = 0.978537048337295
    $ws.Range("M16").Value = 4.823431
    $ws.Range("N16").Value = 14.470293
    $ws.Range("O16").Value = 0.2404505470122564
    $ws.Range("P16").Value = 0.2404505470122564
    $ws.Range("Q16").Value = 21896.39327603696
    $ws.Range("R16").Value = 197067.5394843326
    $ws.Range("S16").Value = 0.2352897685444614
    $ws.Range("T16").Value = 0.2352897685444614
    $ws.Range("G17").Value = 4539.588785666667
    $ws.Range("H17").Value = 13618.766357
    $ws.Range("I17").Value = 0.9785370483372949
    $ws.Range("J17").Value = 0.978537048337295
    $ws.Range("M17").Value = 6.467117666666667
    $ws.Range("N17").Value = 19.401353
    $ws.Range("O17").Value = 0.3223891832479054
    $ws.Range("P17").Value = 0.3223891832479053
    $ws.Range("Q17").Value = 29358.05483518678
    $ws.Range("R17").Value = 264222.493516681
    $ws.Range("S17").Value = 0.3154697597912766
    $ws.Range("T17").Value = 0.3154697597912766
